$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4121
$ws.Range("I62").Value = 3870.5715
$ws.Range("J62").Value = 4997.5
$ws.Range("K62").Value = 3870.5715
$ws.Range("L62").Value = 4997.5
$ws.Range("M62").Value = -3246.5715
$ws.Range("N62").Value = -6245.5
$ws.Range("H64").Value = 3827.5454
$ws.Range("J64").Value = 4017.1667
$ws.Range("L64").Value = 4017.1667
$ws.Range("N64").Value = -4513.1667
$ws.Range("H65").Value = 4121
$ws.Range("I65").Value = 3870.5715
$ws.Range("J65").Value = 4997.5
$ws.Range("K65").Value = 19352.8575
$ws.Range("L65").Value = 24987.5
$ws.Range("M65").Value = -16232.8575
$ws.Range("N65").Value = -31227.5
$ws.Range("H67").Value = 3827.5454
$ws.Range("J67").Value = 4017.1667
$ws.Range("L67").Value = 4017.1667
$ws.Range("N67").Value = -5733.1667
$ws.Range("H70").Value = 13352
$ws.Range("I70").Value = 6146.2
$ws.Range("J70").Value = 18499
$ws.Range("K70").Value = 18438.6
$ws.Range("L70").Value = 55497
$ws.Range("M70").Value = -18168.6
$ws.Range("N70").Value = -56037
$ws.Range("H73").Value = 13352
$ws.Range("I73").Value = 6146.2
$ws.Range("J73").Value = 18499
$ws.Range("K73").Value = 18438.6
$ws.Range("L73").Value = 55497
$ws.Range("M73").Value = -17502.6
$ws.Range("N73").Value = -57369
$ws.Range("H86").Value = 5901.5264
$ws.Range("I86").Value = 6783.875
$ws.Range("J86").Value = 5259.8184
$ws.Range("K86").Value = 6783.875
$ws.Range("L86").Value = 5259.8184
$ws.Range("M86").Value = -5660.875
$ws.Range("N86").Value = -7505.8184
$ws.Range("H89").Value = 5901.5264
$ws.Range("I89").Value = 6783.875
$ws.Range("J89").Value = 5259.8184
$ws.Range("K89").Value = 33919.375
$ws.Range("L89").Value = 26299.092
$ws.Range("M89").Value = -28303.375
$ws.Range("N89").Value = -37531.092
$ws.Range("H116").Value = 5492.647
$ws.Range("I116").Value = 5215.5835
$ws.Range("K116").Value = 5215.5835
$ws.Range("M116").Value = -1773.5835
$ws.Range("H137").Value = 4247.8823
$ws.Range("J137").Value = 1586.4
$ws.Range("L137").Value = 4759.200000000001
$ws.Range("N137").Value = -9859.200000000001
$ws.Range("H138").Value = 2462.83
$ws.Range("I138").Value = 2232.0278
$ws.Range("K138").Value = 6696.0834
$ws.Range("M138").Value = -1556.0834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8519.905000000001
$ws.Range("I32").Value = 7838.0117
$ws.Range("J32").Value = 15035.777
$ws.Range("K32").Value = 7838.0117
$ws.Range("L32").Value = 15035.777
$ws.Range("M32").Value = -7551.0117
$ws.Range("N32").Value = -15609.777
$ws.Range("H61").Value = 3141.2188
$ws.Range("I61").Value = 2212.36
$ws.Range("J61").Value = 6458.5713
$ws.Range("K61").Value = 2212.36
$ws.Range("L61").Value = 6458.5713
$ws.Range("M61").Value = -2000.36
$ws.Range("N61").Value = -6882.5713
$ws.Range("H86").Value = 49999
$ws.Range("J86").Value = 49999
$ws.Range("L86").Value = 49999
$ws.Range("N86").Value = -52371
$ws.Range("H89").Value = 49999
$ws.Range("J89").Value = 49999
$ws.Range("L89").Value = 149997
$ws.Range("N89").Value = -161853
$ws.Range("H136").Value = 3141.2188
$ws.Range("I136").Value = 2212.36
$ws.Range("J136").Value = 6458.5713
$ws.Range("K136").Value = 6637.08
$ws.Range("L136").Value = 19375.7139
$ws.Range("M136").Value = -4087.08
$ws.Range("N136").Value = -24475.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 12822098
$ws.Range("J64").Value = 1493.8
$ws.Range("L64").Value = 1493.8
$ws.Range("N64").Value = -1943.8
$ws.Range("H67").Value = 12822098
$ws.Range("J67").Value = 1493.8
$ws.Range("L67").Value = 1493.8
$ws.Range("N67").Value = -3053.8
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H134").Value = 2191.9119
$ws.Range("I134").Value = 1513.6492
$ws.Range("J134").Value = 5706.5454
$ws.Range("K134").Value = 4540.9476
$ws.Range("L134").Value = 17119.6362
$ws.Range("M134").Value = -2005.9476
$ws.Range("N134").Value = -22189.6362
$ws.Range("H140").Value = 59999.168
$ws.Range("J140").Value = 59999.168
$ws.Range("L140").Value = 59999.168
$ws.Range("N140").Value = -70359.16800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2939.7407
$ws.Range("I31").Value = 1848.75
$ws.Range("K31").Value = 1848.75
$ws.Range("M31").Value = -1553.75
$ws.Range("H34").Value = 2939.7407
$ws.Range("I34").Value = 1848.75
$ws.Range("K34").Value = 1848.75
$ws.Range("M34").Value = -1646.75
$ws.Range("H58").Value = 1463.2285
$ws.Range("I58").Value = 1150.6
$ws.Range("J58").Value = 2244.8
$ws.Range("K58").Value = 1150.6
$ws.Range("L58").Value = 2244.8
$ws.Range("M58").Value = -947.5999999999999
$ws.Range("N58").Value = -2650.8
$ws.Range("H62").Value = 38467908
$ws.Range("I62").Value = 62506210
$ws.Range("J62").Value = 6628.2
$ws.Range("K62").Value = 62506210
$ws.Range("L62").Value = 6628.2
$ws.Range("M62").Value = -62505586
$ws.Range("N62").Value = -7876.2
$ws.Range("H65").Value = 38467908
$ws.Range("I65").Value = 62506210
$ws.Range("J65").Value = 6628.2
$ws.Range("K65").Value = 312531050
$ws.Range("L65").Value = 33141
$ws.Range("M65").Value = -312527930
$ws.Range("N65").Value = -39381
$ws.Range("H125").Value = 52878.668
$ws.Range("J125").Value = 52878.668
$ws.Range("L125").Value = 52878.668
$ws.Range("N125").Value = -57798.668
$ws.Range("H136").Value = 1463.2285
$ws.Range("I136").Value = 1150.6
$ws.Range("J136").Value = 2244.8
$ws.Range("K136").Value = 3451.8
$ws.Range("L136").Value = 6734.400000000001
$ws.Range("M136").Value = -901.7999999999997
$ws.Range("N136").Value = -11834.4
$ws.Range("H141").Value = 275662.12
$ws.Range("I141").Value = 150000
$ws.Range("J141").Value = 291369.88
$ws.Range("K141").Value = 150000
$ws.Range("L141").Value = 291369.88
$ws.Range("M141").Value = -144820
$ws.Range("N141").Value = -301729.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1334354.1
$ws.Range("I9").Value = 5500000
$ws.Range("K9").Value = 16500000
$ws.Range("M9").Value = -16499776
$ws.Range("H131").Value = 2674974
$ws.Range("I131").Value = 4202702
$ws.Range("J131").Value = 1449.75
$ws.Range("K131").Value = 12608106
$ws.Range("L131").Value = 4349.25
$ws.Range("M131").Value = -12603066
$ws.Range("N131").Value = -14429.25
$ws.Range("H132").Value = 2098.5405
$ws.Range("J132").Value = 2336.1785
$ws.Range("L132").Value = 21025.6065
$ws.Range("N132").Value = -26085.6065

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3521.2092
$ws.Range("I102").Value = 2676.7324
$ws.Range("J102").Value = 7518.4
$ws.Range("K102").Value = 2676.7324
$ws.Range("L102").Value = 7518.4
$ws.Range("M102").Value = -1054.7324
$ws.Range("N102").Value = -10762.4
$ws.Range("H122").Value = 2002102.2
$ws.Range("I122").Value = 3335170.2
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 10005510.6
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -10003060.6
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 1968.8605
$ws.Range("I132").Value = 1291.24
$ws.Range("J132").Value = 2910
$ws.Range("K132").Value = 3873.72
$ws.Range("L132").Value = 8730
$ws.Range("M132").Value = -1343.72
$ws.Range("N132").Value = -13790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 33336434
$ws.Range("H71").Value = 33336434
$ws.Range("H74").Value = 46312.332
$ws.Range("I74").Value = 38988.5
$ws.Range("K74").Value = 38988.5
$ws.Range("M74").Value = -37990.5
$ws.Range("H77").Value = 46312.332
$ws.Range("I77").Value = 38988.5
$ws.Range("K77").Value = 116965.5
$ws.Range("M77").Value = -111973.5
$ws.Range("H136").Value = 5811790
$ws.Range("I136").Value = 8577265
$ws.Range("J136").Value = 4292.8
$ws.Range("K136").Value = 25731795
$ws.Range("L136").Value = 12878.4
$ws.Range("M136").Value = -25729245
$ws.Range("N136").Value = -17978.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3337.8262
$ws.Range("I132").Value = 2978.2
$ws.Range("J132").Value = 4482.091
$ws.Range("K132").Value = 8934.599999999999
$ws.Range("L132").Value = 13446.273
$ws.Range("M132").Value = -6404.599999999999
$ws.Range("N132").Value = -18506.273
